$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)
$ws.Activate()

# --- Remove all existing hyperlinks on the sheet (Excel COM only supports
# deleting the whole sheet collection in one go via a range's Hyperlinks) ---
$ws.Range("A1").Hyperlinks.Delete()

# --- Update cell values that changed ---
$ws.Range("C2").Value = "snehasiddela175@gmail.com"
$ws.Range("D2").Value = "S@sneha789"
$ws.Range("E2").Value = "S@sneha789"
$ws.Range("D3").Value = "P@padmu12345"
$ws.Range("D5").Value = "S@sneha1234"
$ws.Range("D6").Value = "S@sneha1234"

# --- Re-create the hyperlinks in their final order/arrangement ---
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:snehasiddela175@gmail.com")
$ws.Range("C2").Font.Color = 16711680
$ws.Range("C2").Font.Underline = 2

$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:padmusham@gmail.com")
$ws.Range("C3").Font.Color = 16711680
$ws.Range("C3").Font.Underline = 2

$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:P@padmu12345")
$ws.Range("D3").Font.Color = 16711680
$ws.Range("D3").Font.Underline = 2

$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:P@padmu")
$ws.Range("E3").Font.Color = 16711680
$ws.Range("E3").Font.Underline = 2

$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:bumisunu@gmail.com")
$ws.Range("C4").Font.Color = 16711680
$ws.Range("C4").Font.Underline = 2

$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:B@bumisuni5143")
$ws.Range("D4").Font.Color = 16711680
$ws.Range("D4").Font.Underline = 2

$ws.Hyperlinks.Add($ws.Range("E4"), "mailto:B@bumisuni5143")
$ws.Range("E4").Font.Color = 16711680
$ws.Range("E4").Font.Underline = 2

$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:S@sneha1234")
$ws.Range("D5").Font.Color = 16711680
$ws.Range("D5").Font.Underline = 2

$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:S@sneha1234")
$ws.Range("D6").Font.Color = 16711680
$ws.Range("D6").Font.Underline = 2

$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:S@sneha789")
$ws.Range("D2").Font.Color = 16711680
$ws.Range("D2").Font.Underline = 2

# --- E2 keeps its hyperlink-like appearance but is no longer a real hyperlink ---
$ws.Range("E2").Font.Color = 16711680
$ws.Range("E2").Font.Underline = 2

# --- C5 is a new, empty, but hyperlink-styled cell ---
$ws.Range("C5").Font.Color = 16711680
$ws.Range("C5").Font.Underline = 2

# --- Selection moves to D5 ---
$ws.Range("D5").Select()
